# Weekly refresh of fruit/vegetable (Achicoria) price records.
# The date (D), volume (J), min/max/avg price (K/L/M), origin (O) and
# Precio $/Kg (P) columns for rows 2-14 are updated to the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44232

$ws.Range("D3").Value = 44251
$ws.Range("J3").Value = 120
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 5000
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 312

$ws.Range("D4").Value = 44208

$ws.Range("D5").Value = 44188
$ws.Range("J5").Value = 210

$ws.Range("D6").Value = 44186

$ws.Range("D7").Value = 44231

$ws.Range("D8").Value = 44292
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 375

$ws.Range("D9").Value = 44230
$ws.Range("J9").Value = 250

$ws.Range("D10").Value = 44215
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 5000
$ws.Range("M10").Value = 5500
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 344

$ws.Range("D11").Value = 44204
$ws.Range("J11").Value = 430
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5500
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 344

$ws.Range("D12").Value = 44189
$ws.Range("J12").Value = 250

$ws.Range("D13").Value = 44210
$ws.Range("J13").Value = 340

$ws.Range("D14").Value = 44187
$ws.Range("J14").Value = 160
